$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set new/changed cell values (A1:C25 region) ---
$ws.Range("B10").Value = "A ciência dos biomateriais é uma atividade multidisciplinar que envolve a medicina, as ciências naturais e as engenharias, delimitando duas grandes áreas: a biotecnologia e a bioengenharia. A disciplina Biomateriais visa prover aos estudantes fundamentos básicos da ciência de biomateriais, dar uma perspectiva sobre os principais biomateriais aplicados em algumas áreas da medicina e contribuir para a compreensão das interações célula-material. Dessa forma, contribuir para o desenvolvimento da área e certamente alavancar a formação de recursos humanos associados a um melhor uso da infra-estrutura já existente."
$ws.Range("C10").Value = "A ciência dos biomateriais é uma atividade multidisciplinar que envolve a medicina, as ciências naturais e as engenharias, delimitando duas grandes áreas: a biotecnologia e a bioengenharia. A disciplina Biomateriais visa prover aos estudantes fundamentos básicos da ciência de biomateriais, dar uma perspectiva sobre os principais biomateriais aplicados em algumas áreas da medicina e contribuir para a compreensão das interações célula-material. Dessa forma, contribuir para o desenvolvimento da área e certamente alavancar a formação de recursos humanos associados a um melhor uso da infra-estrutura já existente."
$ws.Range("A12").Value = "Programa resumido:"
$ws.Range("B12").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("C12").Value = "2166002 - Sandra Giacomin Schneider"
$ws.Range("A13").Value = "Short syllabus:"
$ws.Range("A14").Value = "Programa:"
$ws.Range("B14").Value = "1 - Introdução aos Biomateriais
2 - Interação tecido - implante
3 - Técnicas de modificação de superfície
4 - Técnicas de caracterização biológica
5 - Aspectos práticos no uso de biomateriais"
$ws.Range("C14").Value = "1 - Introdução aos Biomateriais
2 - Interação tecido - implante
3 - Técnicas de modificação de superfície
4 - Técnicas de caracterização biológica
5 - Aspectos práticos no uso de biomateriais"
$ws.Range("A15").Value = "Syllabus:"
$ws.Range("A16").Value = "Avaliação:"
$ws.Range("A17").Value = "Método:"
$ws.Range("B17").Value = "1 - Introdução aos Biomateriais
  1.1- Conceitos básicos de biomateriais; 
  1.2 - Classes de materiais usados na área biomédica;
  1.3 - Classificação dos biomateriais quanto à resposta biológica
2 - Interação tecido  implante:
  2.1 - Histórico da osteointegração; 
  2.2 - Fisiologia do osso;
  2.3 - Natureza da ligação osso-implante;
  2.4 - Aspectos superficiais dos implantes.
3 - Técnicas de modificação da superfície:
  3.1 - Técnicas para criar uma superfície bioativa: cerâmicas bioativas e biovidros, recobrimentos com fosfatos de  cálcio como transportador de proteínas ósseas morfogenéticas;
  3.2 - Técnicas para aumentar a rugosidade superficial: usinagem, ataque ácido, jateamento, aspersão térmica. 
4 - Técnicas de caracterização biológica
  4.1 - Teste em líquido corporal simulado (SBF)
  4.2 - Cultura de células (in vitro)
  4.3  Teste com cobaias (in vivo)
5 - Aspectos práticos no uso de biomateriais
  5.1- Técnicas de esterilização
  5.2  Normas técnicas"
$ws.Range("C17").Value = "1 - Introdução aos Biomateriais
  1.1- Conceitos básicos de biomateriais; 
  1.2 - Classes de materiais usados na área biomédica;
  1.3 - Classificação dos biomateriais quanto à resposta biológica
2 - Interação tecido  implante:
  2.1 - Histórico da osteointegração; 
  2.2 - Fisiologia do osso;
  2.3 - Natureza da ligação osso-implante;
  2.4 - Aspectos superficiais dos implantes.
3 - Técnicas de modificação da superfície:
  3.1 - Técnicas para criar uma superfície bioativa: cerâmicas bioativas e biovidros, recobrimentos com fosfatos de  cálcio como transportador de proteínas ósseas morfogenéticas;
  3.2 - Técnicas para aumentar a rugosidade superficial: usinagem, ataque ácido, jateamento, aspersão térmica. 
4 - Técnicas de caracterização biológica
  4.1 - Teste em líquido corporal simulado (SBF)
  4.2 - Cultura de células (in vitro)
  4.3  Teste com cobaias (in vivo)
5 - Aspectos práticos no uso de biomateriais
  5.1- Técnicas de esterilização
  5.2  Normas técnicas"
$ws.Range("A18").Value = "Critério:"
$ws.Range("B18").Value = "As aulas serão expositivas com auxilio do quadro para anotações e empregando-se recursos audiovisuais."
$ws.Range("C18").Value = "As aulas serão expositivas com auxilio do quadro para anotações e empregando-se recursos audiovisuais."
$ws.Range("A19").Value = "Norma de recuperação:"
$ws.Range("B19").Value = "Serão utilizadas duas notas para compor a nota final sendo: NF=(P1+P2)/2
P1 e P2 serão avaliações escritas (eventualmente a P2 poderá ser substituída por trabalho apresentado por escrito e oral)."
$ws.Range("C19").Value = "Serão utilizadas duas notas para compor a nota final sendo: NF=(P1+P2)/2
P1 e P2 serão avaliações escritas (eventualmente a P2 poderá ser substituída por trabalho apresentado por escrito e oral)."
$ws.Range("A20").Value = "Bibliografia:"
$ws.Range("B20").Value = "Uma prova escrita (Rec) que será composta á NF para obtenção da média final (MF) pelo seguinte critério: 
MF = (Rec+NF)/2"
$ws.Range("C20").Value = "Uma prova escrita (Rec) que será composta á NF para obtenção da média final (MF) pelo seguinte critério: 
MF = (Rec+NF)/2"
$ws.Range("A21").Value = "Requisitos:"
$ws.Range("B22").Value = "LOM3011 -  Ensaios Mecânicos  (Requisito)
"
$ws.Range("C22").Value = "LOM3011 -  Ensaios Mecânicos  (Requisito)
"
$ws.Range("B23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito)
"
$ws.Range("C23").Value = "LOM3013 -  Ciência dos Materiais  (Requisito)
"
$ws.Range("B24").Value = "LOM3036 -  Propriedades Mecânicas  (Requisito)
"
$ws.Range("C24").Value = "LOM3036 -  Propriedades Mecânicas  (Requisito)
"
$ws.Range("B25").Value = "LOM3046 -  Técnicas de Análise Microestrutural  (Requisito)
"
$ws.Range("C25").Value = "LOM3046 -  Técnicas de Análise Microestrutural  (Requisito)
"

# --- Copy formats onto newly-created cells so they pick up the right column style ---
# Column B template style (wrap, normal font) comes from B10; Column C template style (wrap, red font) comes from C10
$ws.Range("B10").Copy() | Out-Null
$ws.Range("B12").PasteSpecial(-4122) | Out-Null
$ws.Range("B14").PasteSpecial(-4122) | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("B22").PasteSpecial(-4122) | Out-Null
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").PasteSpecial(-4122) | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null

# --- Clear cells that must no longer hold content ---
$ws.Range("B13").Clear()
$ws.Range("C13").Clear()
$ws.Range("B15").Clear()
$ws.Range("C15").Clear()
$ws.Range("B21").Clear()
$ws.Range("C21").Clear()
$ws.Range("A22").Clear()

# --- Row heights ---
$ws.Rows.Item(12).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 120
$ws.Rows.Item(16).AutoFit()
$ws.Rows.Item(17).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 120
$ws.Rows.Item(21).AutoFit()
$ws.Rows.Item(22).RowHeight = 30

# --- Remove obsolete row 26 (content now lives in row 25) ---
$ws.Rows.Item(26).Delete()

# --- Column layout: split the old merged A:B width definition so only column A keeps it ---
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(2).ColumnWidth
